$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency values in column C
$ws.Range("C2").Value = 3446
$ws.Range("C3").Value = 3127
$ws.Range("C4").Value = 2127
$ws.Range("C6").Value = 1339
$ws.Range("C7").Value = 699
$ws.Range("C8").Value = 577
$ws.Range("C9").Value = 504
$ws.Range("C10").Value = 491

# Update row 11 labels and value
$ws.Range("A11").Value = "Home Decor"
$ws.Range("B11").Value = "Vintage & Collectibles"
$ws.Range("C11").Value = 456
